$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 146; existing rows 146:214 shift down to 147:215.
$ws.Rows("146:146").Insert()

# Populate the newly inserted row 146 with the new weekly price record.
$ws.Range("A146").Value = 6
$ws.Range("B146").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C146").Value = "Metropolitana"
$ws.Range("D146").Value = 44529
$ws.Range("E146").Value = 13
$ws.Range("F146").Value = 100112026
$ws.Range("G146").Value = "Haba"
$ws.Range("H146").Value = "Sin especificar"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 250
$ws.Range("K146").Value = 7000
$ws.Range("L146").Value = 8000
$ws.Range("M146").Value = 7400
$ws.Range("N146").Value = "$/saco 25 kilos"
$ws.Range("O146").Value = "Región del Maule"
$ws.Range("P146").Value = 296
$ws.Range("Q146").Value = 25
$ws.Range("R146").Value = "Hortaliza"
